$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 472.44446
$ws.Range("I9").Value = 476.125
$ws.Range("K9").Value = 476.125
$ws.Range("M9").Value = -307.125

$ws.Range("H132").Value = 9425
$ws.Range("I132").Value = 4209.5
$ws.Range("K132").Value = 12628.5
$ws.Range("M132").Value = -10098.5

$ws.Range("H135").Value = 2238.2666
$ws.Range("I135").Value = 723
$ws.Range("K135").Value = 6507
$ws.Range("M135").Value = -3972

$ws.Range("H137").Value = 1084710
$ws.Range("I137").Value = 2214.5
$ws.Range("J137").Value = 1301209.1
$ws.Range("K137").Value = 6643.5
$ws.Range("L137").Value = 3903627.3
$ws.Range("M137").Value = -4093.5
$ws.Range("N137").Value = -3908727.3

$ws.Range("H141").Value = 38470540
$ws.Range("I141").Value = 50004820
$ws.Range("J141").Value = 22944
$ws.Range("K141").Value = 150014460
$ws.Range("L141").Value = 68832
$ws.Range("M141").Value = -150009280
$ws.Range("N141").Value = -79192

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4142.035
$ws.Range("I32").Value = 4397.0376
$ws.Range("K32").Value = 4397.0376
$ws.Range("M32").Value = -4110.0376

$ws.Range("H60").Value = 78174
$ws.Range("J60").Value = 24959
$ws.Range("L60").Value = 24959
$ws.Range("N60").Value = -26425

$ws.Range("H61").Value = 5256118
$ws.Range("I61").Value = 5720509.5
$ws.Range("J61").Value = 2005377.8
$ws.Range("K61").Value = 5720509.5
$ws.Range("L61").Value = 2005377.8
$ws.Range("M61").Value = -5720297.5
$ws.Range("N61").Value = -2005801.8

$ws.Range("H74").Value = 545351.4
$ws.Range("I74").Value = 582920.1
$ws.Range("J74").Value = 6865.6665
$ws.Range("K74").Value = 582920.1
$ws.Range("L74").Value = 6865.6665
$ws.Range("M74").Value = -582046.1
$ws.Range("N74").Value = -8613.666499999999

$ws.Range("H77").Value = 545351.4
$ws.Range("I77").Value = 582920.1
$ws.Range("J77").Value = 6865.6665
$ws.Range("K77").Value = 2914600.5
$ws.Range("L77").Value = 34328.3325
$ws.Range("M77").Value = -2910232.5
$ws.Range("N77").Value = -43064.3325

$ws.Range("H110").Value = 8109.923
$ws.Range("I110").Value = 8109.923
$ws.Range("K110").Value = 8109.923
$ws.Range("M110").Value = -6064.923

$ws.Range("H122").Value = 2635.182
$ws.Range("I122").Value = 2158.375
$ws.Range("J122").Value = 3906.6667
$ws.Range("K122").Value = 6475.125
$ws.Range("L122").Value = 11720.0001
$ws.Range("M122").Value = -4025.125
$ws.Range("N122").Value = -16620.0001

$ws.Range("H125").Value = 65000
$ws.Range("J125").Value = 65000
$ws.Range("L125").Value = 65000
$ws.Range("N125").Value = -74840

$ws.Range("H132").Value = 2131972.8
$ws.Range("I132").Value = 4396.725
$ws.Range("J132").Value = 14289550
$ws.Range("K132").Value = 13190.175
$ws.Range("L132").Value = 42868650
$ws.Range("M132").Value = -10660.175
$ws.Range("N132").Value = -42873710

$ws.Range("H136").Value = 5256118
$ws.Range("I136").Value = 5720509.5
$ws.Range("J136").Value = 2005377.8
$ws.Range("K136").Value = 17161528.5
$ws.Range("L136").Value = 6016133.4
$ws.Range("M136").Value = -17158978.5
$ws.Range("N136").Value = -6021233.4

$ws.Range("H141").Value = 135000
$ws.Range("J141").Value = 135000
$ws.Range("L141").Value = 135000
$ws.Range("N141").Value = -145360

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2383.25
$ws.Range("I99").Value = 2383.25
$ws.Range("K99").Value = 2383.25
$ws.Range("M99").Value = -885.25

$ws.Range("H105").Value = 741032.0600000001
$ws.Range("I105").Value = 1091209.4
$ws.Range("J105").Value = 5659.8
$ws.Range("K105").Value = 1091209.4
$ws.Range("L105").Value = 5659.8
$ws.Range("M105").Value = -1089462.4
$ws.Range("N105").Value = -9153.799999999999

$ws.Range("H134").Value = 9093345
$ws.Range("I134").Value = 2644.2222
$ws.Range("K134").Value = 7932.6666
$ws.Range("M134").Value = -5397.6666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 29711718
$ws.Range("I31").Value = 41669316
$ws.Range("K31").Value = 41669316
$ws.Range("M31").Value = -41669021

$ws.Range("H34").Value = 29711718
$ws.Range("I34").Value = 41669316
$ws.Range("K34").Value = 41669316
$ws.Range("M34").Value = -41669114

$ws.Range("H58").Value = 3336.611
$ws.Range("I58").Value = 3216.6428
$ws.Range("K58").Value = 3216.6428
$ws.Range("M58").Value = -3013.6428

$ws.Range("H62").Value = 7359.25
$ws.Range("J62").Value = 8953
$ws.Range("L62").Value = 8953
$ws.Range("N62").Value = -10201

$ws.Range("H65").Value = 7359.25
$ws.Range("J65").Value = 8953
$ws.Range("L65").Value = 44765
$ws.Range("N65").Value = -51005

$ws.Range("H134").Value = 3587
$ws.Range("I134").Value = 3571.2
$ws.Range("J134").Value = 3666
$ws.Range("K134").Value = 10713.6
$ws.Range("L134").Value = 10998
$ws.Range("M134").Value = -8178.599999999999
$ws.Range("N134").Value = -16068

$ws.Range("H136").Value = 3336.611
$ws.Range("I136").Value = 3216.6428
$ws.Range("K136").Value = 9649.928400000001
$ws.Range("M136").Value = -7099.928400000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H97").Value = 1198.9333
$ws.Range("J97").Value = 983.3
$ws.Range("L97").Value = 2949.9
$ws.Range("N97").Value = -3941.9

$ws.Range("H130").Value = 13044.2
$ws.Range("I130").Value = 10040
$ws.Range("K130").Value = 30120
$ws.Range("M130").Value = -25100

$ws.Range("H131").Value = 4918.5264
$ws.Range("I131").Value = 3105.8
$ws.Range("J131").Value = 5565.9287
$ws.Range("K131").Value = 9317.400000000001
$ws.Range("L131").Value = 16697.7861
$ws.Range("M131").Value = -4277.400000000001
$ws.Range("N131").Value = -26777.7861

$ws.Range("H138").Value = 11308.777
$ws.Range("I138").Value = 13055.223
$ws.Range("J138").Value = 9562.333000000001
$ws.Range("K138").Value = 39165.669
$ws.Range("L138").Value = 28686.999
$ws.Range("M138").Value = -34025.669
$ws.Range("N138").Value = -38966.999

$ws.Range("H140").Value = 5899.357
$ws.Range("I140").Value = 2780.6667
$ws.Range("K140").Value = 8342.000100000001
$ws.Range("M140").Value = -3162.000100000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 13622.3
$ws.Range("I70").Value = 11570.667
$ws.Range("J70").Value = 14501.571
$ws.Range("K70").Value = 11570.667
$ws.Range("L70").Value = 14501.571
$ws.Range("M70").Value = -11300.667
$ws.Range("N70").Value = -15041.571

$ws.Range("H73").Value = 13622.3
$ws.Range("I73").Value = 11570.667
$ws.Range("J73").Value = 14501.571
$ws.Range("K73").Value = 11570.667
$ws.Range("L73").Value = 14501.571
$ws.Range("M73").Value = -10634.667
$ws.Range("N73").Value = -16373.571

$ws.Range("H80").Value = 4819.3335
$ws.Range("I80").Value = 3500
$ws.Range("K80").Value = 3500
$ws.Range("M80").Value = -2502

$ws.Range("H83").Value = 4819.3335
$ws.Range("I83").Value = 3500
$ws.Range("K83").Value = 17500
$ws.Range("M83").Value = -12508

$ws.Range("H107").Value = 461.27274
$ws.Range("J107").Value = 457.5
$ws.Range("L107").Value = 457.5
$ws.Range("N107").Value = -4297.5

$ws.Range("H122").Value = 4752.75
$ws.Range("I122").Value = 4707.696
$ws.Range("J122").Value = 4960
$ws.Range("K122").Value = 14123.088
$ws.Range("L122").Value = 14880
$ws.Range("M122").Value = -11673.088
$ws.Range("N122").Value = -19780

$ws.Range("H132").Value = 7694641.5
$ws.Range("I132").Value = 2334
$ws.Range("J132").Value = 33335666
$ws.Range("K132").Value = 7002
$ws.Range("L132").Value = 100006998
$ws.Range("M132").Value = -4472
$ws.Range("N132").Value = -100012058

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H23").Value = 28998.5
$ws.Range("J23").Value = 28998.5
$ws.Range("L23").Value = 28998.5
$ws.Range("N23").Value = -29458.5

$ws.Range("H93").Value = 2344593.8
$ws.Range("I93").Value = 501607.56
$ws.Range("J93").Value = 6952059
$ws.Range("K93").Value = 501607.56
$ws.Range("L93").Value = 6952059
$ws.Range("M93").Value = -500359.56
$ws.Range("N93").Value = -6954555

$ws.Range("H132").Value = 4954.9375
$ws.Range("I132").Value = 2850.3333
$ws.Range("K132").Value = 8550.999899999999
$ws.Range("M132").Value = -6020.999899999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H31").Value = 8000
$ws.Range("J31").Value = 0
$ws.Range("L31").Value = 0
$ws.Range("N31").ClearContents()

$ws.Range("H56").Value = 96642.5
$ws.Range("I56").Value = 128285
$ws.Range("K56").Value = 128285
$ws.Range("M56").Value = -127571

$ws.Range("H81").Value = 1750.3636
$ws.Range("I81").Value = 1706.1111
$ws.Range("J81").Value = 1949.5
$ws.Range("K81").Value = 3412.2222
$ws.Range("L81").Value = 3899
$ws.Range("M81").Value = -2351.2222
$ws.Range("N81").Value = -6021

$ws.Range("H84").Value = 1750.3636
$ws.Range("I84").Value = 1706.1111
$ws.Range("J84").Value = 1949.5
$ws.Range("K84").Value = 17061.111
$ws.Range("L84").Value = 19495
$ws.Range("M84").Value = -11757.111
$ws.Range("N84").Value = -30103

$ws.Range("H132").Value = 250951.7
$ws.Range("I132").Value = 994.51514
$ws.Range("K132").Value = 2983.54542
$ws.Range("M132").Value = -453.5454199999999
